$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = 0.01532826326376212
$ws1.Range("C2").Value = -0.3578207221560517
$ws1.Range("B3").Value = -0.6196767501210488
$ws1.Range("C3").Value = -0.7493574645013042
$ws1.Range("B4").Value = -0.7589033215210698
$ws1.Range("C4").Value = -0.2996017689642435

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.5386628011823984
$ws2.Range("C2").Value = -0.2057924530462218
$ws2.Range("B3").Value = -1.042891470265323
$ws2.Range("C3").Value = -0.005874055797255057
$ws2.Range("B4").Value = -0.7403141285562166
$ws2.Range("C4").Value = 0.6719119038087041
